$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("modelIterations")

# Fill in newly-computed PCA columns (S:Z) for rows 10-15 and 25-27
# (multioutput / PCA-refinement results referenced in the commit message)
$ws.Range("S10").Value = 61.3
$ws.Range("T10").Value = 44.5
$ws.Range("U10").Value = 38.1
$ws.Range("V10").Value = 21.4
$ws.Range("W10").Value = 60.9
$ws.Range("X10").Value = 44.5
$ws.Range("Y10").Value = 37.799999999999997
$ws.Range("Z10").Value = 21.1
$ws.Range("S11").Value = 81
$ws.Range("T11").Value = 61.8
$ws.Range("U11").Value = 55.4
$ws.Range("V11").Value = 30.1
$ws.Range("W11").Value = 80.400000000000006
$ws.Range("X11").Value = 61.9
$ws.Range("Y11").Value = 55.5
$ws.Range("Z11").Value = 30.1
$ws.Range("S12").Value = 0.44
$ws.Range("T12").Value = 0.43
$ws.Range("U12").Value = 0.42
$ws.Range("V12").Value = 0.48
$ws.Range("W12").Value = 0.45
$ws.Range("X12").Value = 0.43
$ws.Range("Y12").Value = 0.42
$ws.Range("Z12").Value = 0.48
$ws.Range("S13").Value = 61.2
$ws.Range("T13").Value = 44.9
$ws.Range("U13").Value = 38.200000000000003
$ws.Range("V13").Value = 21.6
$ws.Range("W13").Value = 61.2
$ws.Range("X13").Value = 44.4
$ws.Range("Y13").Value = 37.799999999999997
$ws.Range("Z13").Value = 20.8
$ws.Range("S14").Value = 80.8
$ws.Range("T14").Value = 62.1
$ws.Range("U14").Value = 55.3
$ws.Range("V14").Value = 30.1
$ws.Range("W14").Value = 80.900000000000006
$ws.Range("X14").Value = 61.6
$ws.Range("Y14").Value = 55.3
$ws.Range("Z14").Value = 30.1
$ws.Range("S15").Value = 0.45
$ws.Range("T15").Value = 0.42
$ws.Range("U15").Value = 0.42
$ws.Range("V15").Value = 0.48
$ws.Range("W15").Value = 0.45
$ws.Range("X15").Value = 0.43
$ws.Range("Y15").Value = 0.42
$ws.Range("Z15").Value = 0.48
$ws.Range("S25").Value = 76
$ws.Range("T25").Value = 54.9
$ws.Range("U25").Value = 50.1
$ws.Range("V25").Value = 29.8
$ws.Range("W25").Value = 75
$ws.Range("X25").Value = 46.6
$ws.Range("Y25").Value = 40.299999999999997
$ws.Range("Z25").Value = 21.7
$ws.Range("S26").Value = 101.3
$ws.Range("T26").Value = 78.5
$ws.Range("U26").Value = 78.7
$ws.Range("V26").Value = 44.2
$ws.Range("W26").Value = 96
$ws.Range("X26").Value = 68.3
$ws.Range("Y26").Value = 59.1
$ws.Range("Z26").Value = 31.5
$ws.Range("S27").Value = 0.13
$ws.Range("T27").Value = 0.07
$ws.Range("U27").Value = -0.17
$ws.Range("V27").Value = -0.12
$ws.Range("W27").Value = 0.22
$ws.Range("X27").Value = 0.3
$ws.Range("Y27").Value = 0.34
$ws.Range("Z27").Value = 0.43

# Remove the scratch "Sheet1" worksheet used for ad-hoc NN results
$wb.Worksheets.Item("Sheet1").Delete()

# Update the view state (scroll position / selection) left behind by the edit session
$ws.Activate()
$win = $excel.ActiveWindow
try { $win.ScrollRow = 7 } catch {}
try { $win.ScrollColumn = 2 } catch {}
try { $win.TopLeftCell = $ws.Range("B7") } catch {}
$ws.Range("AA24").Select()
